$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "fallen",
    "klingen",
    "hauen",
    "sprengen",
    "rasen",
    "trauen",
    "kichern",
    "ehren",
    "biegen",
    "schenken",
    "stecken",
    "gründen",
    "sichern",
    "enden",
    "liegen",
    "wüten",
    "liefern",
    "klettern",
    "warnen",
    "pflanzen",
    "greifen",
    "suchen",
    "mögen",
    "malen",
    "tropfen",
    "stammen",
    "knarren",
    "reizen",
    "fangen",
    "sperren",
    "schmecken",
    "weichen",
    "heilen",
    "schwächen",
    "schwingen",
    "äußern",
    "lügen",
    "bellen",
    "boxen",
    "fahren",
    "rufen",
    "scheinen",
    "bergen",
    "bluten",
    "arten",
    "filmen",
    "münzen",
    "fällen",
    "zünden",
    "kehren",
    "geben",
    "jubeln",
    "dienen",
    "achten",
    "heulen",
    "sinken",
    "zeigen",
    "wehtun",
    "trennen",
    "graben",
    "führen",
    "jagen",
    "treiben",
    "stehlen",
    "ärgern",
    "saufen",
    "betteln",
    "spinnen",
    "runden",
    "altern",
    "spielen",
    "ändern",
    "quälen",
    "platzen",
    "brauchen",
    "kosten",
    "feiern",
    "schrecken",
    "schlucken",
    "spüren",
    "töten",
    "streichen",
    "wundern",
    "zielen",
    "schulden",
    "grüßen",
    "loben",
    "machen",
    "folgen",
    "gelten",
    "freuen",
    "wirken",
    "kümmern",
    "pfeifen",
    "fließen",
    "sorgen",
    "bitten",
    "formen",
    "flüchten",
    "erben",
    "schwören",
    "wachsen",
    "lesen",
    "schreiten",
    "zögern",
    "räumen",
    "wenden",
    "hören",
    "tollen",
    "fischen",
    "werfen",
    "flehen",
    "dringen",
    "bauen",
    "backen",
    "helfen",
    "sterben",
    "drehen",
    "decken",
    "heben",
    "planen",
    "werden",
    "irren",
    "seufzen",
    "lockern",
    "scheitern",
    "siegen",
    "mauern"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $words[$i]
}
